$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44671

$ws.Range("D3").Value = 44965
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 34000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 34600
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("S3").Value = 1922
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44965
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 32000
$ws.Range("O4").Value = 33000
$ws.Range("P4").Value = 32333
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 1796
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44643
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29000
$ws.Range("S5").Value = 1450

$ws.Range("D6").Value = 44679
$ws.Range("M6").Value = 200

$ws.Range("D7").Value = 44679
$ws.Range("L7").Value = 'Tercera'
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("S7").Value = 1225

$ws.Range("D8").Value = 45021
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 22500
$ws.Range("S8").Value = 1125

$ws.Range("D9").Value = 44993
$ws.Range("M9").Value = 130
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25462
$ws.Range("S9").Value = 1273

$ws.Range("D10").Value = 44650
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 31000
$ws.Range("O10").Value = 32000
$ws.Range("P10").Value = 31500
$ws.Range("S10").Value = 1575

$ws.Range("D11").Value = 44650
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 29000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29500
$ws.Range("Q11").Value = '$/caja 20 kilos'
$ws.Range("R11").Value = 'Región de Coquimbo'
$ws.Range("S11").Value = 1475
$ws.Range("T11").Value = 20

$ws.Range("D12").Value = 45028
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("S12").Value = 1075

$ws.Range("D13").Value = 44664
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 29000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 29500
$ws.Range("S13").Value = 1639

$ws.Range("D14").Value = 44636
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 29000
$ws.Range("O14").Value = 30000
$ws.Range("P14").Value = 29500
$ws.Range("Q14").Value = '$/caja 20 kilos'
$ws.Range("S14").Value = 1475
$ws.Range("T14").Value = 20

$ws.Range("D15").Value = 45007
$ws.Range("L15").Value = 'Segunda'
$ws.Range("N15").Value = 27000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 27500
$ws.Range("S15").Value = 1375

$ws.Range("D16").Value = 44972
$ws.Range("M16").Value = 140
$ws.Range("N16").Value = 27000
$ws.Range("O16").Value = 28000
$ws.Range("P16").Value = 27429
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1524
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 45014
$ws.Range("L17").Value = 'Segunda'
$ws.Range("N17").Value = 24000
$ws.Range("O17").Value = 25000
$ws.Range("P17").Value = 24500
$ws.Range("S17").Value = 1225

$ws.Range("D18").Value = 44979
$ws.Range("M18").Value = 250
$ws.Range("Q18").Value = '$/caja 20 kilos'
$ws.Range("S18").Value = 1475
$ws.Range("T18").Value = 20

